# Auto-generated edit script to update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on price cells whose values would otherwise be
# re-interpreted as numbers by Excel (losing exact text formatting, e.g. trailing zeros)
$textCells = @("D5","D6","D7","D8","D9","D10","D11","D15","D18","D19","D20","D21","D23","D24","D25","D26","D28","D31","D32","D33","D34","D36","D37","D38","D39","D41","D42","D44","D45","D46","D47","D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "57.884.58"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "2.551.87"
$ws.Range("E3").Value = "  -3.59%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "518.87"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "138.85"
$ws.Range("E6").Value = "  -3.66%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "0.562"
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("D9").Value = "6.52"
$ws.Range("E9").Value = "  -6.59%  "
$ws.Range("D10").Value = "0.0990"
$ws.Range("E10").Value = "  -3.40%  "
$ws.Range("D11").Value = "0.324"
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("D13").Value = "3.000.80"
$ws.Range("E13").Value = "  -3.55%  "
$ws.Range("D14").Value = "57.850.41"
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("D15").Value = "19.97"
$ws.Range("E15").Value = "  -4.90%  "
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("D17").Value = "2.542.36"
$ws.Range("E17").Value = "  -4.15%  "
$ws.Range("D18").Value = "333.74"
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("D19").Value = "4.29"
$ws.Range("E19").Value = "  -1.75%  "
$ws.Range("D20").Value = "10.12"
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("D21").Value = "6.12"
$ws.Range("E21").Value = "  -3.72%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "64.91"
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("D24").Value = "0.164"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "0.401"
$ws.Range("E25").Value = "  -4.11%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "2.685.67"
$ws.Range("E27").Value = "  -3.16%  "
$ws.Range("D28").Value = "6.94"
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("D29").Value = "0.0₃0754"
$ws.Range("E29").Value = "  -5.71%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "6.17"
$ws.Range("E31").Value = "  -7.39%  "
$ws.Range("D32").Value = "1.57"
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("D33").Value = "149.07"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "18.47"
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("E35").Value = "  -4.26%  "
$ws.Range("D36").Value = "1.13"
$ws.Range("E36").Value = "  -5.20%  "
$ws.Range("D37").Value = "0.832"
$ws.Range("E37").Value = "  -6.52%  "
$ws.Range("D38").Value = "35.71"
$ws.Range("E38").Value = "  -2.70%  "
$ws.Range("D39").Value = "0.820"
$ws.Range("E39").Value = "  -5.53%  "
$ws.Range("E40").Value = "  -4.64%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "3.46"
$ws.Range("E42").Value = "  -3.32%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "0.0954"
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("D45").Value = "0.579"
$ws.Range("E45").Value = "  -5.99%  "
$ws.Range("D46").Value = "259.46"
$ws.Range("E46").Value = "  -5.62%  "
$ws.Range("D47").Value = "0.0519"
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "18.53"
$ws.Range("E48").Value = "  -6.90%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "1.981.55"
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("E50").Value = "  -3.07%  "
$ws.Range("E51").Value = "  -5.14%  "
